$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("VanityFairV2")
$wsData.Select()

$wsData.Columns.Item(1).Hidden = $false
$wsData.Columns.Item(2).Hidden = $false
$wsData.Columns.Item(1).ColumnWidth = 14.6640625
$wsData.Columns.Item(2).ColumnWidth = 18.44140625
